$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the Efficiency-related measurement columns (D3:I3) ---
$ws.Range("D3").Value = 0.13387018000000001
$ws.Range("E3").Value = 1.9611641
$ws.Range("F3").Value = 0.0665002
$ws.Range("G3").Value = 1.3930377
$ws.Range("H3").Value = 0.04991426
$ws.Range("I3").Value = 1.0804151

# Running time value was revised
$ws.Range("J3").Value = 2.13

# --- Row 4: fill in the Efficiency-related measurement columns (D4:I4) ---
$ws.Range("D4").Value = 0.10365949000000001
$ws.Range("E4").Value = 1.5478578999999999
$ws.Range("F4").Value = 0.062288925
$ws.Range("G4").Value = 1.2127775000000001
$ws.Range("H4").Value = 0.04759176
$ws.Range("I4").Value = 1.102904

# --- New "Efficiency" column header ---
$ws.Range("K1").Value = "Efficiency"

# --- New "Efficiency" formulas: Efficiency = 1 / ((Fixed + Variable) * RunningTime) ---
$ws.Range("K3").Formula = "=1/((H3+I3)*J3)"
$ws.Range("K4:K8").Formula = "=1/((H4+I4)*J4)"

# --- Update the active selection to the newly-added efficiency cell ---
$ws.Range("K3").Select()
